# Szenario2 weighting sheet: remove the "barrierefreiheit" criterion row
# from the "multifunktionale_nutzungsqualitaet" group, letting the
# remaining rows re-flow and re-balance to equal thirds (1/3) weights
# within their groups, matching the updated R-script (MCA Baeume) logic.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gewichtung")

# Remove the "barrierefreiheit" / "Befahrbarkeit"-group data row (row 10:
# group=multifunktionale_nutzungsqualitaet, criterion=barrierefreiheit).
# Shifting the remaining rows up moves the former rows 11-13
# (kreislauffaehigkeit group) into rows 10-12.
$ws.Range("A10:G10").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)

# multifunktionale_nutzungsqualitaet now has 3 criteria left (rows 7-9):
# recompute within_group_weight as equal thirds.
$ws.Range("F7").Formula = "=1/3"
$ws.Range("F8:F9").Formula = "=1/3"

# kreislauffaehigkeit group (now rows 10-12) keeps equal-thirds weighting.
$ws.Range("F10").Formula = "=1/3"
$ws.Range("F11:F12").Formula = "=1/3"

# Restore the active selection to match the saved view state.
$ws.Range("D18").Select() | Out-Null
